# Fruta / hortaliza, semanal
# A new weekly observation for Femacal de La Calera / Ciboulette is inserted
# at the top of the date-ordered data block (row 278), pushing the existing
# rows 278-381 down to 279-382 and extending the used range to A1:R382.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 278; Excel shifts rows 278:381 down to 279:382
# and carries the formatting of the row above (including the date style
# on column D) onto the newly inserted cells.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(278, 1).Value = 3
$ws.Cells.Item(278, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(278, 3).Value = "Coquimbo"
$ws.Cells.Item(278, 4).Value = 44900
$ws.Cells.Item(278, 5).Value = 5
$ws.Cells.Item(278, 6).Value = 100112039
$ws.Cells.Item(278, 7).Value = "Ciboulette"
$ws.Cells.Item(278, 8).Value = "Sin especificar"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 120
$ws.Cells.Item(278, 11).Value = 1500
$ws.Cells.Item(278, 12).Value = 1500
$ws.Cells.Item(278, 13).Value = 1500
$ws.Cells.Item(278, 14).Value = "$/docena de atados"
$ws.Cells.Item(278, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(278, 16).Value = 500
$ws.Cells.Item(278, 17).Value = 3
$ws.Cells.Item(278, 18).Value = "Hortaliza"
